$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.606831789016724
$ws.Range("B1").Value = 1.837117433547974
$ws.Range("C1").Value = 5.114257335662842
$ws.Range("D1").Value = 1.894317626953125
$ws.Range("E1").Value = 0.6346949338912964
